# Update Hgf-St14 NATMI TPM output with the newly regenerated data.
# The "ECs" sending-cluster rows are gone (only FAPs and MuSCs now send
# signal); target clusters stay ECs / FAPs / MuSCs for each remaining
# sender, so the table shrinks from 9 data rows to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing three rows (old MuSCs-sender block) that no longer exist.
$ws.Rows("8:10").Delete()

# New 6-row body (rows 2-7). Columns A-T as in the sheet header.
$rows = @(
    @("FAPs","Hgf","St14","ECs",3,1,8.583520999999999,25.750563,0.8910607110509009,0.8910607110509009,3,1,0.1404576666666667,0.421373,0.183302244487095,0.1833022444870949,1.205621331444333,10.850591982999,0.1633334283098969,0.1633334283098969),
    @("FAPs","Hgf","St14","FAPs",3,1,8.583520999999999,25.750563,0.8910607110509009,0.8910607110509009,3,1,0.5312866666666666,1.59386,0.6933479729318232,0.6933479729318232,4.560310260353332,41.04279234318,0.6178151377663311,0.6178151377663311),
    @("FAPs","Hgf","St14","MuSCs",3,1,8.583520999999999,25.750563,0.8910607110509009,0.8910607110509009,2,0.6666666666666666,0.09451833333333333,0.283555,0.1233497825810819,0.1233497825810819,0.8113000990516666,7.301700891465,0.1099121449746728,0.1099121449746728),
    @("MuSCs","Hgf","St14","ECs",3,1,1.049404,3.148212,0.1089392889490991,0.1089392889490991,3,1,0.1404576666666667,0.421373,0.183302244487095,0.1833022444870949,0.1473968372306667,1.326571535076,0.01996881617719804,0.01996881617719803),
    @("MuSCs","Hgf","St14","FAPs",3,1,1.049404,3.148212,0.1089392889490991,0.1089392889490991,3,1,0.5312866666666666,1.59386,0.6933479729318232,0.6933479729318232,0.5575343531466666,5.017809178319999,0.075532835165492,0.075532835165492),
    @("MuSCs","Hgf","St14","MuSCs",3,1,1.049404,3.148212,0.1089392889490991,0.1089392889490991,2,0.6666666666666666,0.09451833333333333,0.283555,0.1233497825810819,0.1233497825810819,0.09918791707333333,0.89269125366,0.01343763760640902,0.01343763760640902)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
